# Weekly update: a new, more recent price record is inserted at the top of
# the data (row 11), pushing the existing rows 11-31 down to 12-32. The
# previously-last row keeps all of its data (it simply ends up at row 32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 11; this shifts rows 11-31 down to 12-32
# and extends the sheet dimension from R31 to R32 automatically.
$ws.Rows("11").Insert()

# Populate the newly inserted row 11 with the new weekly record.
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = 'Vega Monumental Concepción'
$ws.Range("C11").Value = 'Bíobío'
$ws.Range("D11").Value = 44540
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 100112031
$ws.Range("G11").Value = 'Poroto verde'
$ws.Range("H11").Value = 'Magnum'
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 170
$ws.Range("K11").Value = 21000
$ws.Range("L11").Value = 22000
$ws.Range("M11").Value = 21529
$ws.Range("N11").Value = '$/saco 25 kilos'
$ws.Range("O11").Value = 'Región Metropolitana'
$ws.Range("P11").Value = 861
$ws.Range("Q11").Value = 25
$ws.Range("R11").Value = 'Hortaliza'

# Make sure the new date cell carries the same date/time number format used
# by the rest of column D (style index 2 in the original workbook).
$ws.Range("D11").NumberFormat = $ws.Range("D12").NumberFormat
